$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: copy cell formats onto the positions that need a newly-formatted cell ---
# (done first, while the source cells still hold their original content/format)
$ws.Range("B14").Copy() | Out-Null
$ws.Range("B12").PasteSpecial(-4122) | Out-Null
$ws.Range("C14").Copy() | Out-Null
$ws.Range("C12").PasteSpecial(-4122) | Out-Null
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A13").PasteSpecial(-4122) | Out-Null
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A14").PasteSpecial(-4122) | Out-Null
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A15").PasteSpecial(-4122) | Out-Null
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A16").PasteSpecial(-4122) | Out-Null
$ws.Range("B22").Copy() | Out-Null
$ws.Range("B18").PasteSpecial(-4122) | Out-Null
$ws.Range("C22").Copy() | Out-Null
$ws.Range("C18").PasteSpecial(-4122) | Out-Null
$ws.Range("B22").Copy() | Out-Null
$ws.Range("B20").PasteSpecial(-4122) | Out-Null
$ws.Range("C22").Copy() | Out-Null
$ws.Range("C20").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Step 2: clear cells that must end up empty (their old content is no longer used) ---
$ws.Range("B13").Clear() | Out-Null
$ws.Range("C13").Clear() | Out-Null
$ws.Range("B15").Clear() | Out-Null
$ws.Range("C15").Clear() | Out-Null
$ws.Range("B16").Clear() | Out-Null
$ws.Range("C16").Clear() | Out-Null
$ws.Range("A22").Clear() | Out-Null

# --- Step 3: write the new cell values (rows 10-22 final content) ---
$ws.Range("A10").Value = "Objetivos:"
$ws.Range("B10").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C10").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("A11").Value = "Objectives:"
$ws.Range("A12").Value = "Programa resumido:"
$ws.Range("B12").Value = "5840963 - Daniela Camargo Vernilli"
$ws.Range("C12").Value = "5840963 - Daniela Camargo Vernilli"
$ws.Range("A13").Value = "Short syllabus:"
$ws.Range("A14").Value = "Programa:"
$ws.Range("B14").Value = "6495737 - Durval Rodrigues Junior"
$ws.Range("C14").Value = "6495737 - Durval Rodrigues Junior"
$ws.Range("A15").Value = "Syllabus:"
$ws.Range("A16").Value = "Avaliação:"
$ws.Range("A17").Value = "Método:"
$ws.Range("B17").Value = "984972 - Hugo Ricardo Zschommler Sandim"
$ws.Range("C17").Value = "984972 - Hugo Ricardo Zschommler Sandim"
$ws.Range("A18").Value = "Critério:"
$ws.Range("B18").Value = "Aulas expositivas complementadas com experimentos desenvolvidos em laboratório didático; realização de relatórios para cada experimento e de estudo de casos."
$ws.Range("C18").Value = "Aulas expositivas complementadas com experimentos desenvolvidos em laboratório didático; realização de relatórios para cada experimento e de estudo de casos."
$ws.Range("A19").Value = "Norma de recuperação:"
$ws.Range("B19").Value = "Média aritmética das notas obtidas nos relatórios e trabalhos. Será aprovado o aluno que obtiver nota final maior ou igual a 5,0."
$ws.Range("C19").Value = "Média aritmética das notas obtidas nos relatórios e trabalhos. Será aprovado o aluno que obtiver nota final maior ou igual a 5,0."
$ws.Range("A20").Value = "Bibliografia:"
$ws.Range("B20").Value = "Devido às características práticas da disciplina, não será oferecida recuperação."
$ws.Range("C20").Value = "Devido às características práticas da disciplina, não será oferecida recuperação."
$ws.Range("A21").Value = "Requisitos:"
$ws.Range("B22").Value = "LOB1012 -  Estatística  (Requisito fraco)`n"
$ws.Range("C22").Value = "LOB1012 -  Estatística  (Requisito fraco)`n"

# --- Step 4: remove the now-unused trailing rows (23-27) ---
$ws.Range("A23:A27").EntireRow.Delete() | Out-Null

# --- Step 5: fix up row heights on the restructured rows ---
$ws.Rows.Item(12).RowHeight = 60
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(14).RowHeight = 120
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(17).RowHeight = 60
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(19).RowHeight = 60
$ws.Rows.Item(20).RowHeight = 120
$ws.Rows.Item(22).RowHeight = 30
